$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("B3:F3").Value = ""
$ws.Range("H13").Select() | Out-Null
